$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D ("Test Steps") gets new example/placeholder text, italic style ---
# Apply italic font to the whole column (header + data, D1:D9) first so the
# new style ends up as cellXfs index 2 (matching fontId=2 / <i/>).
$ws.Range("D1:D9").Font.Italic = $true

$ws.Range("D2").Value = "<InvestigationName>,<InvestigationDescription>"
$ws.Range("D3").Value = "<InvestigationName>"
$ws.Range("D4").Value = "<InvestigationNameOld>,<InvestigationNameNew><InvestigationDescription>"
$ws.Range("D5").Value = "<InvestigationName>,<SetName>"
$ws.Range("D6").Value = "<InvestigationName>,<ExplortionName>"
$ws.Range("D7").Value = "<InvestigationName>,<Comparison>"
$ws.Range("D8").Value = "<InvestigationName>,<tobeShared>"
$ws.Range("D9").Value = "<InvestigationName>,<toBeDeleted>"

# --- Column E ("Test Data") updated/added sample data ---
$ws.Range("E3").Value = "Gene_invet"
$ws.Range("E4").Value = "Investigation1,modified, To analysis the impact of the virus"
$ws.Range("E5").Value = "InvSet,Set1"
$ws.Range("E6").Value = "InvExploration,Explo1"
$ws.Range("E7").Value = "InvComparison,Comp1"
$ws.Range("E8").Value = "InvShare,divya.devanathan@eaglegenomics.com"
$ws.Range("E9").Value = "InvDelete"

# --- Column widths: C narrower, D much wider (now holds the longer placeholders) ---
# (Target stored widths are 21.42578125 / 50.5703125; the engine quantizes
# ColumnWidth to a 1/6-character grid, so these inputs are the closest
# achievable values.)
$ws.Columns.Item(3).ColumnWidth = 20.666666666666668
$ws.Columns.Item(4).ColumnWidth = 49.666666666666664

# --- Selection / view: active cell moves from E7 to D6 ---
[void]$ws.Range("D6").Select()
